$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add the new data row (Sheet1!A35:E35) ---
$ws.Range("A35").Value = 290
$ws.Range("B35").Value = 4
$ws.Range("C35").Value = 469
$ws.Range("D35").Value = 6275
$ws.Range("E35").Formula = "=C35/D35"

# --- Extend the two line-chart series so they pick up the new row ---
$co1 = $ws.ChartObjects().Item(1)
$chart1 = $co1.Chart
$series1 = $chart1.SeriesCollection().Item(1)
$series1.Formula = "=SERIES(Sheet1!`$B`$1,Sheet1!`$A`$2:`$A`$35,Sheet1!`$B`$2:`$B`$35,1)"

$co2 = $ws.ChartObjects().Item(2)
$chart2 = $co2.Chart
$series2 = $chart2.SeriesCollection().Item(1)
$series2.Formula = "=SERIES(Sheet1!`$E`$1,Sheet1!`$A`$2:`$A`$35,Sheet1!`$E`$2:`$E`$35,1)"

# --- Match the updated view / selection state ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("U14").Select() | Out-Null
